# Guru99_NewCustomerCreationPage_FRAMEWORK.xlsx
# "minor changes fetching dd values from excel"
#
# - AutoCompleteSampleSheet!A2 test-data value "java" -> "JavaScript"
# - That sheet becomes the active/selected tab (was DataFromSeleniumEasyURL),
#   with the active cell moving to A10.

$wb = $excel.ActiveWorkbook

$autoComplete = $wb.Worksheets.Item("AutoCompleteSampleSheet")
$autoComplete.Range("A2").Value = "JavaScript"

$autoComplete.Activate()
$autoComplete.Range("A10").Select()
